# Add "Wins" / "Losses" / "Ties" season-record columns (AD:AF) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - bold/bordered/centered header style matches the
# existing header cells (e.g. A1) so copy that style across.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Determine the last data row (season record is the same for every player
# row: the team's overall record for the season).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 92  # AD
    $ws.Cells.Item($r, 31).Value = 70  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
